$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1581
$ws.Range("I40").Value = 1424.0646
$ws.Range("J40").Value = 1928.5
$ws.Range("K40").Value = 1424.0646
$ws.Range("L40").Value = 1928.5
$ws.Range("M40").Value = -1249.0646
$ws.Range("N40").Value = -2278.5

$ws.Range("H57").Value = 41695.25
$ws.Range("J57").Value = 41695.25
$ws.Range("L57").Value = 125085.75
$ws.Range("N57").Value = -126083.75

$ws.Range("H86").Value = 2676.7058
$ws.Range("I86").Value = 1508.3334
$ws.Range("K86").Value = 1508.3334
$ws.Range("M86").Value = -385.3334

$ws.Range("H89").Value = 2676.7058
$ws.Range("I89").Value = 1508.3334
$ws.Range("K89").Value = 7541.666999999999
$ws.Range("M89").Value = -1925.666999999999

$ws.Range("H113").Value = 2645
$ws.Range("I113").Value = 2700
$ws.Range("K113").Value = 2700
$ws.Range("M113").Value = 554

$ws.Range("H132").Value = 3655
$ws.Range("I132").Value = 3786
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 11358
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -8828
$ws.Range("N132").Value = -14060

$ws.Range("H137").Value = 1333.5758
$ws.Range("I137").Value = 1174.64
$ws.Range("J137").Value = 1830.25
$ws.Range("K137").Value = 3523.92
$ws.Range("L137").Value = 5490.75
$ws.Range("M137").Value = -973.9200000000001
$ws.Range("N137").Value = -10590.75

$ws.Range("H141").Value = 1959.4762
$ws.Range("I141").Value = 1907.45
$ws.Range("K141").Value = 5722.35
$ws.Range("M141").Value = -542.3500000000004

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3358.1558
$ws.Range("I32").Value = 2588.5232
$ws.Range("J32").Value = 7527
$ws.Range("K32").Value = 2588.5232
$ws.Range("L32").Value = 7527
$ws.Range("M32").Value = -2301.5232
$ws.Range("N32").Value = -8101

$ws.Range("H45").Value = 10112.333
$ws.Range("I45").Value = 10984.8
$ws.Range("J45").Value = 5750
$ws.Range("K45").Value = 10984.8
$ws.Range("L45").Value = 5750
$ws.Range("M45").Value = -10607.8
$ws.Range("N45").Value = -6504

$ws.Range("H61").Value = 11695.182
$ws.Range("I61").Value = 12664.7
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 12664.7
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -12452.7
$ws.Range("N61").Value = -2424

$ws.Range("H74").Value = 2745.8667
$ws.Range("I74").Value = 2711
$ws.Range("J74").Value = 2785.7144
$ws.Range("K74").Value = 2711
$ws.Range("L74").Value = 2785.7144
$ws.Range("M74").Value = -1837
$ws.Range("N74").Value = -4533.7144

$ws.Range("H77").Value = 2745.8667
$ws.Range("I77").Value = 2711
$ws.Range("J77").Value = 2785.7144
$ws.Range("K77").Value = 13555
$ws.Range("L77").Value = 13928.572
$ws.Range("M77").Value = -9187
$ws.Range("N77").Value = -22664.572

$ws.Range("H102").Value = 5292692.5
$ws.Range("I102").Value = 6174307.5
$ws.Range("J102").Value = 3000
$ws.Range("K102").Value = 6174307.5
$ws.Range("L102").Value = 3000
$ws.Range("M102").Value = -6172685.5
$ws.Range("N102").Value = -6244

$ws.Range("H122").Value = 5129587
$ws.Range("I122").Value = 5129587
$ws.Range("K122").Value = 15388761
$ws.Range("M122").Value = -15386311

$ws.Range("H132").Value = 4226.7827
$ws.Range("I132").Value = 1927.7273
$ws.Range("J132").Value = 6334.25
$ws.Range("K132").Value = 5783.1819
$ws.Range("L132").Value = 19002.75
$ws.Range("M132").Value = -3253.1819
$ws.Range("N132").Value = -24062.75

$ws.Range("H136").Value = 11695.182
$ws.Range("I136").Value = 12664.7
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 37994.10000000001
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -35444.10000000001
$ws.Range("N136").Value = -11100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 10754426
$ws.Range("I86").Value = 15153063
$ws.Range("K86").Value = 15153063
$ws.Range("M86").Value = -15151940

$ws.Range("H89").Value = 10754426
$ws.Range("I89").Value = 15153063
$ws.Range("K89").Value = 75765315
$ws.Range("M89").Value = -75759699

$ws.Range("H107").Value = 1082
$ws.Range("I107").Value = 1150.25
$ws.Range("J107").Value = 900
$ws.Range("K107").Value = 1150.25
$ws.Range("L107").Value = 900
$ws.Range("M107").Value = 769.75
$ws.Range("N107").Value = -4740

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 581.6667
$ws.Range("I22").Value = 497.5
$ws.Range("J22").Value = 750
$ws.Range("K22").Value = 497.5
$ws.Range("L22").Value = 750
$ws.Range("M22").Value = -147.5
$ws.Range("N22").Value = -1450

$ws.Range("H25").Value = 200000940
$ws.Range("I25").Value = 1177.75
$ws.Range("J25").Value = 1000000000
$ws.Range("K25").Value = 1177.75
$ws.Range("L25").Value = 1000000000
$ws.Range("M25").Value = -1003.75
$ws.Range("N25").Value = -1000000348

$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()

$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()

$ws.Range("H105").Value = 1086.0667
$ws.Range("I105").Value = 842
$ws.Range("K105").Value = 842
$ws.Range("M105").Value = 905

$ws.Range("H132").Value = 2416.5527
$ws.Range("I132").Value = 2331.9375
$ws.Range("J132").Value = 2867.8333
$ws.Range("K132").Value = 6995.8125
$ws.Range("L132").Value = 8603.499899999999
$ws.Range("M132").Value = -4465.8125
$ws.Range("N132").Value = -13663.4999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H43").Value = 7224.75
$ws.Range("J43").Value = 8966.666999999999
$ws.Range("L43").Value = 26900.001
$ws.Range("N43").Value = -27128.001

$ws.Range("H118").Value = 1648.6428
$ws.Range("I118").Value = 545.2857
$ws.Range("J118").Value = 2752
$ws.Range("K118").Value = 1635.8571
$ws.Range("L118").Value = 8256
$ws.Range("M118").Value = -392.8571000000002
$ws.Range("N118").Value = -10742

$ws.Range("H132").Value = 4044.5833
$ws.Range("I132").Value = 986.6667
$ws.Range("J132").Value = 4481.4287
$ws.Range("K132").Value = 8880.0003
$ws.Range("L132").Value = 40332.85830000001
$ws.Range("M132").Value = -6350.0003
$ws.Range("N132").Value = -45392.85830000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1331.5333
$ws.Range("I102").Value = 1138.875
$ws.Range("J102").Value = 1551.7142
$ws.Range("K102").Value = 1138.875
$ws.Range("L102").Value = 1551.7142
$ws.Range("M102").Value = 483.125
$ws.Range("N102").Value = -4795.7142

$ws.Range("H113").Value = 71429800
$ws.Range("I113").Value = 166667440
$ws.Range("J113").Value = 1562.875
$ws.Range("K113").Value = 166667440
$ws.Range("L113").Value = 1562.875
$ws.Range("M113").Value = -166665270
$ws.Range("N113").Value = -5902.875

$ws.Range("H122").Value = 2494704.5
$ws.Range("I122").Value = 2702450.8
$ws.Range("J122").Value = 1750
$ws.Range("K122").Value = 8107352.399999999
$ws.Range("L122").Value = 5250
$ws.Range("M122").Value = -8104902.399999999
$ws.Range("N122").Value = -10150

$ws.Range("H126").Value = 8634.200000000001
$ws.Range("I126").Value = 11781.2
$ws.Range("J126").Value = 2340.2
$ws.Range("K126").Value = 35343.60000000001
$ws.Range("L126").Value = 7020.599999999999
$ws.Range("M126").Value = -32873.60000000001
$ws.Range("N126").Value = -11960.6

$ws.Range("H132").Value = 4637.5835
$ws.Range("I132").Value = 8666
$ws.Range("J132").Value = 2978.8235
$ws.Range("K132").Value = 25998
$ws.Range("L132").Value = 8936.470499999999
$ws.Range("M132").Value = -23468
$ws.Range("N132").Value = -13996.4705

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2709.9
$ws.Range("I7").Value = 2466.3333
$ws.Range("J7").Value = 2814.2856
$ws.Range("K7").Value = 2466.3333
$ws.Range("L7").Value = 2814.2856
$ws.Range("M7").Value = -2354.3333
$ws.Range("N7").Value = -3038.2856

$ws.Range("H40").Value = 66668828
$ws.Range("I40").Value = 90910960
$ws.Range("K40").Value = 90910960
$ws.Range("M40").Value = -90910824

$ws.Range("H61").Value = 3067.8572
$ws.Range("I61").Value = 2700
$ws.Range("J61").Value = 3558.3333
$ws.Range("K61").Value = 2700
$ws.Range("L61").Value = 3558.3333
$ws.Range("M61").Value = -2498
$ws.Range("N61").Value = -3962.3333

$ws.Range("H100").Value = 667.55554
$ws.Range("I100").Value = 515.4286
$ws.Range("K100").Value = 515.4286
$ws.Range("M100").Value = 25.57140000000004

$ws.Range("H113").Value = 3067.8572
$ws.Range("I113").Value = 2700
$ws.Range("J113").Value = 3558.3333
$ws.Range("K113").Value = 2700
$ws.Range("L113").Value = 3558.3333
$ws.Range("M113").Value = -530
$ws.Range("N113").Value = -7898.3333

$ws.Range("H122").Value = 20358342
$ws.Range("J122").Value = 10000000
$ws.Range("L122").Value = 30000000
$ws.Range("N122").Value = -30004900

$ws.Range("H126").Value = 2709.9
$ws.Range("I126").Value = 2466.3333
$ws.Range("J126").Value = 2814.2856
$ws.Range("K126").Value = 7398.999899999999
$ws.Range("L126").Value = 8442.856800000001
$ws.Range("M126").Value = -4928.999899999999
$ws.Range("N126").Value = -13382.8568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 90909750
$ws.Range("I107").Value = 200000500
$ws.Range("J107").Value = 803.8333
$ws.Range("K107").Value = 600001500
$ws.Range("L107").Value = 2411.4999
$ws.Range("M107").Value = -599999580
$ws.Range("N107").Value = -6251.4999

$ws.Range("H135").Value = 48000
$ws.Range("J135").Value = 48000
$ws.Range("L135").Value = 48000
$ws.Range("N135").Value = -58140

$ws.Range("H136").Value = 2849.0334
$ws.Range("I136").Value = 3457.5557
$ws.Range("K136").Value = 10372.6671
$ws.Range("M136").Value = -7822.667099999999
